$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: insert a new data row at row 5 (pushing existing rows 5-13
# down to 6-14) and populate it with the latest week's observation, carrying
# forward the market/product metadata that is constant across every row.
$ws.Rows("5:5").Insert()

$ws.Range("A5").Value2 = 8
$ws.Range("B5").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C5").Value2 = "Coquimbo"
$ws.Range("D5").Value2 = 44484
$ws.Range("E5").Value2 = 4
$ws.Range("F5").Value2 = 100112026
$ws.Range("G5").Value2 = "Haba"
$ws.Range("H5").Value2 = "Sin especificar"
$ws.Range("I5").Value2 = "Primera"
$ws.Range("J5").Value2 = 400
$ws.Range("K5").Value2 = 9000
$ws.Range("L5").Value2 = 10000
$ws.Range("M5").Value2 = 9500
$ws.Range("N5").Value2 = "`$/saco 25 kilos"
$ws.Range("O5").Value2 = "Provincia del Elquí"
$ws.Range("P5").Value2 = 380
$ws.Range("Q5").Value2 = 25
$ws.Range("R5").Value2 = "Hortaliza"
